# Apply updated "dSF" (column F) values as per repull/mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 0
    4  = 2
    5  = 4
    6  = 1
    10 = 6
    11 = 4
    13 = 9
    15 = -4
    20 = 3
    22 = -3
    23 = 7
    24 = -4
    27 = -3
    31 = 1
    33 = -2
    40 = 3
    45 = 5
    46 = -3
    47 = -5
    52 = -2
    53 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
